# Update LR-pair statistics for Adam9-Itga6 sheet per Dr Hou advice
# Ligand/Receptor-expressing cell counts changed from 1 to 3, and the
# dependent expression/specificity/edge-weight metrics were recomputed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.32689766666667
$ws.Range("H2").Value = 30.980693
$ws.Range("I2").Value = 0.2044815006034941
$ws.Range("J2").Value = 0.204481500603494
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 147.4213356666667
$ws.Range("N2").Value = 442.264007
$ws.Range("O2").Value = 0.9507885170992249
$ws.Range("P2").Value = 0.950788517099225
$ws.Range("Q2").Value = 1522.405047312983
$ws.Range("R2").Value = 13701.64542581685
$ws.Range("S2").Value = 0.1944186627330204
$ws.Range("T2").Value = 0.1944186627330204
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.32689766666667
$ws.Range("H3").Value = 30.980693
$ws.Range("I3").Value = 0.2044815006034941
$ws.Range("J3").Value = 0.204481500603494
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.340788333333334
$ws.Range("N3").Value = 7.022365000000001
$ws.Range("O3").Value = 0.01509682881537204
$ws.Range("P3").Value = 0.01509682881537204
$ws.Range("Q3").Value = 24.17308157766056
$ws.Range("R3").Value = 217.557734198945
$ws.Range("S3").Value = 0.003087022210521343
$ws.Range("T3").Value = 0.003087022210521343
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.32689766666667
$ws.Range("H4").Value = 30.980693
$ws.Range("I4").Value = 0.2044815006034941
$ws.Range("J4").Value = 0.204481500603494
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.289533666666667
$ws.Range("N4").Value = 15.868601
$ws.Range("O4").Value = 0.03411465408540306
$ws.Range("P4").Value = 0.03411465408540307
$ws.Range("Q4").Value = 54.62447288005478
$ws.Range("R4").Value = 491.620255920493
$ws.Range("S4").Value = 0.006975815659952337
$ws.Range("T4").Value = 0.006975815659952338
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.37031933333333
$ws.Range("H5").Value = 88.110958
$ws.Range("I5").Value = 0.5815577111671272
$ws.Range("J5").Value = 0.5815577111671272
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 147.4213356666667
$ws.Range("N5").Value = 442.264007
$ws.Range("O5").Value = 0.9507885170992249
$ws.Range("P5").Value = 0.950788517099225
$ws.Range("Q5").Value = 4329.811705076522
$ws.Range("R5").Value = 38968.30534568871
$ws.Range("S5").Value = 0.5529383938082122
$ws.Range("T5").Value = 0.5529383938082122
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 29.37031933333333
$ws.Range("H6").Value = 88.110958
$ws.Range("I6").Value = 0.5815577111671272
$ws.Range("J6").Value = 0.5815577111671272
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.340788333333334
$ws.Range("N6").Value = 7.022365000000001
$ws.Range("O6").Value = 0.01509682881537204
$ws.Range("P6").Value = 0.01509682881537204
$ws.Range("Q6").Value = 68.74970084174112
$ws.Range("R6").Value = 618.74730757567
$ws.Range("S6").Value = 0.008779677211749692
$ws.Range("T6").Value = 0.008779677211749692
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 29.37031933333333
$ws.Range("H7").Value = 88.110958
$ws.Range("I7").Value = 0.5815577111671272
$ws.Range("J7").Value = 0.5815577111671272
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.289533666666667
$ws.Range("N7").Value = 15.868601
$ws.Range("O7").Value = 0.03411465408540306
$ws.Range("P7").Value = 0.03411465408540307
$ws.Range("Q7").Value = 155.3552929144176
$ws.Range("R7").Value = 1398.197636229758
$ws.Range("S7").Value = 0.01983964014716529
$ws.Range("T7").Value = 0.01983964014716529
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.80562866666667
$ws.Range("H8").Value = 32.416886
$ws.Range("I8").Value = 0.2139607882293788
$ws.Range("J8").Value = 0.2139607882293788
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 147.4213356666667
$ws.Range("N8").Value = 442.264007
$ws.Range("O8").Value = 0.9507885170992249
$ws.Range("P8").Value = 0.950788517099225
$ws.Range("Q8").Value = 1592.980210758022
$ws.Range("R8").Value = 14336.8218968222
$ws.Range("S8").Value = 0.2034314605579924
$ws.Range("T8").Value = 0.2034314605579924
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.80562866666667
$ws.Range("H9").Value = 32.416886
$ws.Range("I9").Value = 0.2139607882293788
$ws.Range("J9").Value = 0.2139607882293788
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.340788333333334
$ws.Range("N9").Value = 7.022365000000001
$ws.Range("O9").Value = 0.01509682881537204
$ws.Range("P9").Value = 0.01509682881537204
$ws.Range("Q9").Value = 25.29368951726556
$ws.Range("R9").Value = 227.64320565539
$ws.Range("S9").Value = 0.003230129393101
$ws.Range("T9").Value = 0.003230129393100999
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 10.80562866666667
$ws.Range("H10").Value = 32.416886
$ws.Range("I10").Value = 0.2139607882293788
$ws.Range("J10").Value = 0.2139607882293788
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.289533666666667
$ws.Range("N10").Value = 15.868601
$ws.Range("O10").Value = 0.03411465408540306
$ws.Range("P10").Value = 0.03411465408540307
$ws.Range("Q10").Value = 57.15673662183178
$ws.Range("R10").Value = 514.410629596486
$ws.Range("S10").Value = 0.007299198278285437
$ws.Range("T10").Value = 0.007299198278285437
